# Fruta / hortaliza, semanal
# New weekly record added for "Vega Modelo de Temuco" - Mango.
# A new row is inserted at row 170 (pushing the existing rows 170-255 down
# to 171-256, with the last existing row re-appearing as the new row 256),
# and the freshly inserted row 170 is populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 170; this shifts rows 170:255
# down to 171:256 (carrying all their values/formatting with them).
$ws.Rows("170:170").Insert()

# Make sure the date cell keeps the same date/time number format used by
# every other "Fecha" cell in column D.
$ws.Cells.Item(170, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the newly inserted row 170 with this week's data.
$ws.Cells.Item(170, 1).Value = 10                                      # A - Mercado ID
$ws.Cells.Item(170, 2).Value = "Vega Modelo de Temuco"                 # B - Mercado
$ws.Cells.Item(170, 3).Value = "La Araucanía"                          # C - Región
$ws.Cells.Item(170, 4).Value = 44529                                   # D - Fecha
$ws.Cells.Item(170, 5).Value = 9                                       # E - Codreg
$ws.Cells.Item(170, 6).Value = "Fruta"                                 # F - Tipo
$ws.Cells.Item(170, 7).Value = 100108                                  # G - Producto ID
$ws.Cells.Item(170, 8).Value = "Tropicales y subtropicales"            # H - Producto
$ws.Cells.Item(170, 9).Value = 100108002                               # I - Categoría ID
$ws.Cells.Item(170, 10).Value = "Mango"                                # J - Categoría
$ws.Cells.Item(170, 11).Value = "Sin especificar"                      # K - Variedad
$ws.Cells.Item(170, 12).Value = "Primera"                              # L - Calidad
$ws.Cells.Item(170, 13).Value = 450                                    # M - Volumen
$ws.Cells.Item(170, 14).Value = 7000                                   # N - Precio mínimo
$ws.Cells.Item(170, 15).Value = 7000                                   # O - Precio máximo
$ws.Cells.Item(170, 16).Value = 7000                                   # P - Precio promedio ponderado
$ws.Cells.Item(170, 17).Value = "$/bandeja 4 kilos"                    # Q - Unidad de comercialización
$ws.Cells.Item(170, 18).Value = "Perú"                                 # R - Origen
$ws.Cells.Item(170, 19).Value = 1750                                   # S - Precio $/Kg
$ws.Cells.Item(170, 20).Value = 4                                      # T - Kg / unidad
